# Sync attendance_reports: swap the "Recorded By" (column G) name order
# from "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# for every row where that exact value occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7th column
    $val = $cell.Text
    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}
